$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2070.8572
$ws.Range("I4").Value = 416
$ws.Range("K4").Value = 416
$ws.Range("M4").Value = -302

$ws.Range("H29").Value = 1602.5
$ws.Range("J29").Value = 2737.5
$ws.Range("L29").Value = 8212.5
$ws.Range("N29").Value = -8774.5

$ws.Range("H38").Value = 2540.9473
$ws.Range("J38").Value = 4995.5557
$ws.Range("L38").Value = 14986.6671
$ws.Range("N38").Value = -15730.6671

$ws.Range("H58").Value = 250022.17
$ws.Range("I58").Value = 44.666668
$ws.Range("J58").Value = 499999.66
$ws.Range("K58").Value = 134.000004
$ws.Range("L58").Value = 1499998.98
$ws.Range("M58").Value = 15.99999600000001
$ws.Range("N58").Value = -1500298.98

$ws.Range("H80").Value = 5605.8335
$ws.Range("I80").Value = 984.5454999999999
$ws.Range("J80").Value = 9516.154
$ws.Range("K80").Value = 2953.6365
$ws.Range("L80").Value = 28548.462
$ws.Range("M80").Value = -1955.6365
$ws.Range("N80").Value = -30544.462

$ws.Range("H83").Value = 5605.8335
$ws.Range("I83").Value = 984.5454999999999
$ws.Range("J83").Value = 9516.154
$ws.Range("K83").Value = 8860.9095
$ws.Range("L83").Value = 85645.386
$ws.Range("M83").Value = -3868.9095
$ws.Range("N83").Value = -95629.386

$ws.Range("H86").Value = 304500.7
$ws.Range("I86").Value = 203800.6
$ws.Range("J86").Value = 405200.8
$ws.Range("K86").Value = 203800.6
$ws.Range("L86").Value = 405200.8
$ws.Range("M86").Value = -202677.6
$ws.Range("N86").Value = -407446.8

$ws.Range("H89").Value = 304500.7
$ws.Range("I89").Value = 203800.6
$ws.Range("J89").Value = 405200.8
$ws.Range("K89").Value = 1019003
$ws.Range("L89").Value = 2026004
$ws.Range("M89").Value = -1013387
$ws.Range("N89").Value = -2037236

$ws.Range("H138").Value = 1545953.8
$ws.Range("I138").Value = 3374.4614
$ws.Range("J138").Value = 1981900.1
$ws.Range("K138").Value = 10123.3842
$ws.Range("L138").Value = 5945700.300000001
$ws.Range("M138").Value = -4983.3842
$ws.Range("N138").Value = -5955980.300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 30669.666
$ws.Range("J9").Value = 30669.666
$ws.Range("L9").Value = 30669.666
$ws.Range("N9").Value = -31009.666

$ws.Range("H20").Value = 30669.666
$ws.Range("J20").Value = 30669.666
$ws.Range("L20").Value = 30669.666
$ws.Range("N20").Value = -31209.666

$ws.Range("H25").Value = 6766.8
$ws.Range("I25").Value = 5954
$ws.Range("J25").Value = 10018
$ws.Range("K25").Value = 5954
$ws.Range("L25").Value = 10018
$ws.Range("M25").Value = -5552
$ws.Range("N25").Value = -10822

$ws.Range("H32").Value = 21215.555
$ws.Range("I32").Value = 14260
$ws.Range("J32").Value = 39300
$ws.Range("K32").Value = 14260
$ws.Range("L32").Value = 39300
$ws.Range("M32").Value = -13973
$ws.Range("N32").Value = -39874

$ws.Range("H122").Value = 8930945
$ws.Range("I122").Value = 2128.625
$ws.Range("J122").Value = 20836034
$ws.Range("K122").Value = 6385.875
$ws.Range("L122").Value = 62508102
$ws.Range("M122").Value = -3935.875
$ws.Range("N122").Value = -62513002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 18221.428
$ws.Range("I134").Value = 2005.88
$ws.Range("K134").Value = 6017.64
$ws.Range("M134").Value = -3482.64

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5327.3228
$ws.Range("I31").Value = 2132
$ws.Range("J31").Value = 5941.8076
$ws.Range("K31").Value = 2132
$ws.Range("L31").Value = 5941.8076
$ws.Range("M31").Value = -1837
$ws.Range("N31").Value = -6531.8076

$ws.Range("H34").Value = 5327.3228
$ws.Range("I34").Value = 2132
$ws.Range("J34").Value = 5941.8076
$ws.Range("K34").Value = 2132
$ws.Range("L34").Value = 5941.8076
$ws.Range("M34").Value = -1930
$ws.Range("N34").Value = -6345.8076

$ws.Range("H99").Value = 2209.0908
$ws.Range("I99").Value = 2230
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 2230
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -732
$ws.Range("N99").Value = -4996

$ws.Range("H126").Value = 2209.0908
$ws.Range("I126").Value = 2230
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 6690
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -4220
$ws.Range("N126").Value = -10940

$ws.Range("H132").Value = 2679.9565
$ws.Range("I132").Value = 2324.5
$ws.Range("J132").Value = 3959.6
$ws.Range("K132").Value = 6973.5
$ws.Range("L132").Value = 11878.8
$ws.Range("M132").Value = -4443.5
$ws.Range("N132").Value = -16938.8

$ws.Range("H134").Value = 3466.8909
$ws.Range("I134").Value = 2428.3914
$ws.Range("J134").Value = 4213.3125
$ws.Range("K134").Value = 7285.174199999999
$ws.Range("L134").Value = 12639.9375
$ws.Range("M134").Value = -4750.174199999999
$ws.Range("N134").Value = -17709.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 23843524
$ws.Range("J9").Value = 23843524
$ws.Range("L9").Value = 71530572
$ws.Range("N9").Value = -71531020

$ws.Range("H10").Value = 750
$ws.Range("I10").Value = 562.5
$ws.Range("J10").Value = 1500
$ws.Range("K10").Value = 1687.5
$ws.Range("L10").Value = 4500
$ws.Range("M10").Value = -1548.5
$ws.Range("N10").Value = -4778

$ws.Range("H22").Value = 142858400
$ws.Range("I22").Value = 250000800
$ws.Range("J22").Value = 1866.6666
$ws.Range("K22").Value = 750002400
$ws.Range("L22").Value = 5599.9998
$ws.Range("M22").Value = -750002231
$ws.Range("N22").Value = -5937.9998

$ws.Range("H27").Value = 142858400
$ws.Range("I27").Value = 250000800
$ws.Range("J27").Value = 1866.6666
$ws.Range("K27").Value = 750002400
$ws.Range("L27").Value = 5599.9998
$ws.Range("M27").Value = -750002298
$ws.Range("N27").Value = -5803.9998

$ws.Range("H61").Value = 139.64285
$ws.Range("I61").Value = 76
$ws.Range("J61").Value = 298.75
$ws.Range("K61").Value = 228
$ws.Range("L61").Value = 896.25
$ws.Range("M61").Value = -13
$ws.Range("N61").Value = -1326.25

$ws.Range("H114").Value = 895.6
$ws.Range("I114").Value = 826
$ws.Range("K114").Value = 2478
$ws.Range("M114").Value = 776

$ws.Range("H120").Value = 6444.727
$ws.Range("I120").Value = 7843.3335
$ws.Range("J120").Value = 5920.25
$ws.Range("K120").Value = 23530.0005
$ws.Range("L120").Value = 17760.75
$ws.Range("M120").Value = -18692.0005
$ws.Range("N120").Value = -27436.75

$ws.Range("H131").Value = 598.42
$ws.Range("I131").Value = 288.1579
$ws.Range("J131").Value = 1009.6977
$ws.Range("K131").Value = 864.4737
$ws.Range("L131").Value = 3029.0931
$ws.Range("M131").Value = 4175.5263
$ws.Range("N131").Value = -13109.0931

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

$ws.Range("H59").Value = 35000
$ws.Range("J59").Value = 35000
$ws.Range("L59").Value = 35000
$ws.Range("N59").Value = -36166

$ws.Range("H122").Value = 39999.5
$ws.Range("I122").Value = 50000
$ws.Range("J122").Value = 29999
$ws.Range("K122").Value = 150000
$ws.Range("L122").Value = 89997
$ws.Range("M122").Value = -147550
$ws.Range("N122").Value = -94897

$ws.Range("H132").Value = 5594.5454
$ws.Range("I132").Value = 3981.1428
$ws.Range("J132").Value = 6347.467
$ws.Range("K132").Value = 11943.4284
$ws.Range("L132").Value = 19042.401
$ws.Range("M132").Value = -9413.428400000001
$ws.Range("N132").Value = -24102.401

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3570.0908
$ws.Range("I122").Value = 1343.8
$ws.Range("K122").Value = 4031.4
$ws.Range("M122").Value = -1581.4
